# Regen save_data: recompute column G ("K", formerly "Strike#") per row.
# The new K values replace the previously-saved ones after the std/mean
# regen + s_vals calc step described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K (column G) value.
$kValues = @{
    2  = 3
    3  = 0
    4  = 0
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 0
    14 = 1
    15 = 0
    16 = 0
    17 = 0
    18 = 1
    19 = 2
    20 = 2
    21 = 0
    22 = 3
    23 = 1
    24 = 3
    25 = 1
    26 = 1
    27 = 1
    28 = 1
    29 = 1
    30 = 1
    31 = 1
    32 = 0
    33 = 0
    34 = 1
    35 = 0
    36 = 0
    37 = 0
    38 = 2
    39 = 2
    40 = 2
    41 = 3
    42 = 2
    43 = 1
    44 = 1
    45 = 1
    46 = 2
    47 = 0
    48 = 1
    49 = 0
    50 = 1
    51 = 2
    52 = 1
    53 = 2
    54 = 2
    55 = 1
    56 = 0
    57 = 1
    58 = 3
    59 = 2
    60 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
